# Weekly update: insert a new price record for "Piña" at row 346 of the
# "Macroferia Regional de Talca" sheet. All existing rows from 346 onward
# shift down by one (346->347, ..., 371->372); the new row takes the data
# for the latest market day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 346:371 down to 347:372, carrying formatting/number-format along.
$ws.Rows.Item(346).Insert()

# Populate the newly-opened row 346 with this week's entry.
$ws.Cells.Item(346, 1).Value  = 5
$ws.Cells.Item(346, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(346, 3).Value  = "Maule"
$ws.Cells.Item(346, 4).Value  = "04/05/2023"
$ws.Cells.Item(346, 5).Value  = 7
$ws.Cells.Item(346, 6).Value  = "Fruta"
$ws.Cells.Item(346, 7).Value  = 100108
$ws.Cells.Item(346, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(346, 9).Value  = 100108005
$ws.Cells.Item(346, 10).Value = "Piña"
$ws.Cells.Item(346, 11).Value = "Caramelo"
$ws.Cells.Item(346, 12).Value = "Segunda"
$ws.Cells.Item(346, 13).Value = 230
$ws.Cells.Item(346, 14).Value = 19000
$ws.Cells.Item(346, 15).Value = 19000
$ws.Cells.Item(346, 16).Value = 19000
$ws.Cells.Item(346, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(346, 18).Value = "Ecuador"
$ws.Cells.Item(346, 19).Value = 1357
$ws.Cells.Item(346, 20).Value = 14
